$d = $word.ActiveDocument

function Split-IntoDefinitionTerm([string]$anchorText, [string]$termText) {
    # Finds "<anchorText><termText>" where termText is the trailing term
    # (e.g. " Email"), removes the leading space, breaks the paragraph
    # there, and restyles the new trailing paragraph as DefinitionTerm.
    $r = $d.Content
    $found = $r.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
    if (-not $found) {
        throw "Could not find anchor text: $anchorText"
    }
    $r.Collapse(0)
    $found2 = $r.Find.Execute(" " + $termText, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
    if (-not $found2) {
        throw "Could not find term text: $termText"
    }

    $spaceRange = $d.Range($r.Start, $r.Start + 1)
    $spaceRange.Delete()

    $breakPoint = $d.Range($r.Start, $r.Start)
    $breakPoint.InsertAfter("`r")

    $termPara = $breakPoint.Paragraphs(1).Next()
    $termPara.Style = "Definition Term"
}

# 1. "...Check out the course Moodle here." [space]"Email" -> split
Split-IntoDefinitionTerm "Check out the course Moodle here." "Email"

# 2. "...regularly." [space]"Slack" -> split
Split-IntoDefinitionTerm "you must check your adelphi.edu email account regularly." "Slack"

# 3. "...one of the instructors." [space]"Podcast manager" -> split
Split-IntoDefinitionTerm "one of the instructors." "Podcast manager"

# 4. soundcloud RSS link [space]"Padlet" -> split
Split-IntoDefinitionTerm "sounds.rss" "Padlet"
